# Add "2022-Q4" worksheet (fund-holdings detail, same shape as the other
# quarterly sheets) positioned before "2022-Q3", and insert the matching
# summary row at the top of "总计".
#
# xlPasteFormats
$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (totals) sheet: insert a new row 2 for "2022-Q4" and push the
#    existing "2022-Q3"/"2022-Q2" rows down by one.
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# capture the two existing data rows (read with Value2 - Value's getter
# is unreliable in this host) before we start overwriting cells.
$oldR2B = $tot.Range("B2").Value2
$oldR2C = $tot.Range("C2").Value2
$oldR2D = $tot.Range("D2").Value2
$oldR3B = $tot.Range("B3").Value2
$oldR3C = $tot.Range("C3").Value2
$oldR3D = $tot.Range("D3").Value2

# row4 is brand new - clone row3's formatting onto it first (column A
# carries the bordered/centered "index" style) then fill in the old
# row3 ("2022-Q2") values, shifted down one row.
$tot.Range("A3").Copy()
$tot.Range("A4").PasteSpecial($xlPasteFormats)
$tot.Cells.Item(4, 1).Value = 2
$tot.Cells.Item(4, 2).Value = $oldR3B
$tot.Cells.Item(4, 3).Value = $oldR3C
$tot.Cells.Item(4, 4).Value = $oldR3D

# row3 becomes the old row2 ("2022-Q3") data.
$tot.Cells.Item(3, 1).Value = 1
$tot.Cells.Item(3, 2).Value = $oldR2B
$tot.Cells.Item(3, 3).Value = $oldR2C
$tot.Cells.Item(3, 4).Value = $oldR2D

# row2 becomes the new "2022-Q4" summary entry.
$tot.Cells.Item(2, 1).Value = 0
$tot.Cells.Item(2, 2).Value = "2022-Q4"
$tot.Cells.Item(2, 3).Value = 5
$tot.Cells.Item(2, 4).Value = 0.08

# ---------------------------------------------------------------------
# 2) New "2022-Q4" sheet: duplicate "2022-Q3" (same columns/formatting)
#    right before it, rename it, then overwrite the changed cells with
#    the Q4 fund-holdings figures.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

function Set-TextValue($ws, $addr, $text) {
    # Force a literal-text cell (keeps leading zeros / avoids numeric
    # auto-coercion of numeric-looking strings like "1.26" or "002288"),
    # then strip the quote-prefix style stamp PasteSpecial leaves behind
    # so the cell's formatting matches its already-unstyled neighbours.
    $ws.Range($addr).Value = "'" + $text
    $ws.Range("B2").Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

# Row 2 (008602 / 方正富邦新兴成长混合A) — only D:H change.
Set-TextValue $q4 "D2" "1.26"
Set-TextValue $q4 "E2" "87.26"
Set-TextValue $q4 "F2" "3.44"
Set-TextValue $q4 "G2" "0.0433"
$q4.Range("H2").Value = 8

# Row 3 — now 163823 / 中银稳健策略灵活配置混合.
Set-TextValue $q4 "B3" "163823"
$q4.Range("C3").Value = "中银稳健策略灵活配置混合"
Set-TextValue $q4 "D3" "2.02"
Set-TextValue $q4 "E3" "31.34"
Set-TextValue $q4 "F3" "1.44"
Set-TextValue $q4 "G3" "0.0291"
$q4.Range("H3").Value = 6

# Row 4 — now 002288 / 中银稳进策略灵活配置混合A.
Set-TextValue $q4 "B4" "002288"
$q4.Range("C4").Value = "中银稳进策略灵活配置混合A"
Set-TextValue $q4 "D4" "0.61"
Set-TextValue $q4 "E4" "23.98"
Set-TextValue $q4 "F4" "0.74"
Set-TextValue $q4 "G4" "0.0045"
$q4.Range("H4").Value = 10

# Row 5 (008603 / 方正富邦新兴成长混合C) — only E:H change.
Set-TextValue $q4 "E5" "87.26"
Set-TextValue $q4 "F5" "3.44"
Set-TextValue $q4 "G5" "0.0010"
$q4.Range("H5").Value = 8

# Row 6 — now 016520 / 中银稳进策略灵活配置混合C.
Set-TextValue $q4 "B6" "016520"
$q4.Range("C6").Value = "中银稳进策略灵活配置混合C"
Set-TextValue $q4 "D6" "0.00"
Set-TextValue $q4 "E6" "23.98"
Set-TextValue $q4 "F6" "0.74"
$q4.Range("G6").Value = 0
$q4.Range("H6").Value = 10
